$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated s_val data (filtered save games) for row 2
$ws.Range("B2").Value = 3.230985683306322
$ws.Range("C2").Value = 10.29869402782916
$ws.Range("D2").Value = 26.21740644021617
$ws.Range("E2").Value = 9353990175.932438
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 9353990215.679523
